# Update NATMI LR-pair output with newly computed TPM-based values.
# The "ECs" cluster label is renamed to "Resolving-Mac", and the per-row
# Target cluster assignments / receptor & edge statistics are refreshed
# to match the new TPM run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending cluster = FAPs, Target cluster -> FAPs)
$ws.Range("D2").Value = "FAPs"
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.987076
$ws.Range("N2").Value = 11.961228
$ws.Range("O2").Value = 0.2813308272685638
$ws.Range("P2").Value = 0.2813308272685638
$ws.Range("Q2").Value = 0.2243354891906667
$ws.Range("R2").Value = 2.019019402716
$ws.Range("S2").Value = 0.07676470769453013
$ws.Range("T2").Value = 0.07676470769453013

# Row 3 (Sending cluster = FAPs, Target cluster -> MuSCs)
$ws.Range("D3").Value = "MuSCs"
$ws.Range("M3").Value = 10.131229
$ws.Range("N3").Value = 30.393687
$ws.Range("O3").Value = 0.7148664925918803
$ws.Range("P3").Value = 0.7148664925918804
$ws.Range("Q3").Value = 0.5700403538376666
$ws.Range("R3").Value = 5.130363184539
$ws.Range("S3").Value = 0.1950604485019465
$ws.Range("T3").Value = 0.1950604485019465

# Row 4 (Sending cluster = FAPs, Target cluster -> Resolving-Mac)
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05389233333333333
$ws.Range("N4").Value = 0.161677
$ws.Range("O4").Value = 0.00380268013955587
$ws.Range("P4").Value = 0.00380268013955587
$ws.Range("Q4").Value = 0.003032288063222222
$ws.Range("R4").Value = 0.027290592569
$ws.Range("S4").Value = 0.001037609821159545
$ws.Range("T4").Value = 0.001037609821159545

# Row 5 (Sending cluster = MuSCs, Target cluster -> FAPs)
$ws.Range("D5").Value = "FAPs"
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.987076
$ws.Range("N5").Value = 11.961228
$ws.Range("O5").Value = 0.2813308272685638
$ws.Range("P5").Value = 0.2813308272685638
$ws.Range("Q5").Value = 0.5978195173893334
$ws.Range("R5").Value = 5.380375656504
$ws.Range("S5").Value = 0.2045661195740336
$ws.Range("T5").Value = 0.2045661195740336

# Row 6 (Sending cluster = MuSCs, Target cluster -> MuSCs)
$ws.Range("D6").Value = "MuSCs"
$ws.Range("M6").Value = 10.131229
$ws.Range("N6").Value = 30.393687
$ws.Range("O6").Value = 0.7148664925918803
$ws.Range("P6").Value = 0.7148664925918804
$ws.Range("Q6").Value = 1.519069722107333
$ws.Range("R6").Value = 13.671627498966
$ws.Range("S6").Value = 0.5198060440899338
$ws.Range("T6").Value = 0.5198060440899339

# Row 7 (Sending cluster = MuSCs, Target cluster -> Resolving-Mac)
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.05389233333333333
$ws.Range("N7").Value = 0.161677
$ws.Range("O7").Value = 0.00380268013955587
$ws.Range("P7").Value = 0.00380268013955587
$ws.Range("Q7").Value = 0.008080580531777777
$ws.Range("R7").Value = 0.07272522478600001
$ws.Range("S7").Value = 0.002765070318396325
$ws.Range("T7").Value = 0.002765070318396325

Write-Host "Updated Lama1-Itgb8 LR-pair sheet with new TPM values"
